{"js": "// \"added space after link\" \u2014 insert a space right after the Udacity hyperlink\n// text (outside the hyperlink run, underline turned off so the space itself\n// isn't underlined), leaving the rest of the line (\"| Frontend, Backend,\n// Flask, CRUD, deployment\") untouched.\n\n// 1. Locate the hyperlink text. Word's search matches the logical text even\n//    though it is currently split across three runs (\"Udacity Full Stack Web\n//    Develop\" + \"e\" + \"r Nano Degree\"); inserting right after it will also\n//    cause those runs to coalesce into a single run, same as Word does when\n//    you type at that location and re-save.\nconst linkResults = context.document.body.search(\n  \"Udacity Full Stack Web Developer Nano Degree\",\n  { matchCase: true }\n);\nlinkResults.load(\"items\");\nawait context.sync();\n\nif (linkResults.items.length === 0) {\n  throw new Error(\"Could not find the Udacity hyperlink text\");\n}\n\nconst linkRange = linkResults.items[0];\n\n// Insert a single space immediately after the hyperlink text.\nconst spaceRange = linkRange.insertText(\" \", Word.InsertLocation.after);\n// The new run inherits the hyperlink's underline; explicitly clear it so the\n// space is not underlined (matches the \"<w:u w:val=\"none\"/>\" in the target).\nspaceRange.font.underline = Word.UnderlineType.none;\nawait context.sync();\n\n// 2. Word keeps a single \"_GoBack\" bookmark that marks the location of the\n//    most recent edit. Replicate that bookkeeping: remove it from wherever it\n//    currently sits and drop it at the new edit position (right after \"Fro\"\n//    within \"Frontend\", which is where the cursor would be once the sentence\n//    is reflowed/saved).\nconst goBack = context.document.getBookmarkRangeOrNullObject(\"_GoBack\");\ngoBack.load(\"isNullObject\");\nawait context.sync();\nif (!goBack.isNullObject) {\n  context.document.deleteBookmark(\"_GoBack\");\n}\n\nconst froResults = context.document.body.search(\"Fro\", { matchCase: true });\nfroResults.load(\"items\");\nawait context.sync();\n\nif (froResults.items.length > 0) {\n  const froRange = froResults.items[0];\n  const afterFro = froRange.getRange(Word.RangeLocation.after);\n  afterFro.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# \"added space after link\" \u2014 insert a space right after the Udacity\n# hyperlink text (outside the hyperlink, with underline turned off so the\n# space itself isn't underlined), leaving the rest of the line\n# (\"| Frontend, Backend, Flask, CRUD, deployment\") untouched.\n\n$d = $word.ActiveDocument\n\n# 1. Locate the hyperlink text. Find matches the logical text even though it\n#    is currently split across three runs (\"Udacity Full Stack Web Develop\" +\n#    \"e\" + \"r Nano Degree\"); inserting right after it also coalesces those\n#    runs into a single run, same as Word does when you type at that spot and\n#    re-save.\n$rng = $d.Content\n$found = $rng.Find.Execute(\"Udacity Full Stack Web Developer Nano Degree\")\nif (-not $found) {\n    throw \"Could not find the Udacity hyperlink text\"\n}\n\n# InsertAfter leaves $rng pointing at the original matched text (it does not\n# grow to include the inserted text), so $rng.End is exactly the boundary\n# between the hyperlink and the newly inserted character.\n$rng.InsertAfter(\" \")\n$spaceRange = $d.Range($rng.End, $rng.End + 1)\n\n# The new run inherits the hyperlink's underline; explicitly clear it so the\n# space is not underlined (matches the \"<w:u w:val=\"none\"/>\" in the target).\n$spaceRange.Font.Underline = 0   # wdUnderlineNone\n\n# 2. Word keeps a single \"_GoBack\" bookmark that marks the location of the\n#    most recent edit. Replicate that bookkeeping: remove it from wherever it\n#    currently sits and drop it at the new edit position (right after \"Fro\"\n#    within \"Frontend\", which is where the cursor would be once the sentence\n#    is reflowed/saved).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n$froRng = $d.Content\n$froFound = $froRng.Find.Execute(\"Fro\")\nif ($froFound) {\n    $newBookmarkRange = $d.Range($froRng.End, $froRng.End)\n    $d.Bookmarks.Add(\"_GoBack\", $newBookmarkRange)\n}\n"}
